# Update result data for Zeeland (re-run of optimization with new result files)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: updated source .h5 result file paths (EmissionLimit Brownfield columns E/F/G)
$ws.Range("E5").Value = "Z:\AdOpt_NET0\AdOpt_results\MY\EmissionLimit Brownfield\Zeeland\20250428094955_2030_minC_DD10-1\optimization_results.h5"
$ws.Range("F5").Value = "Z:\AdOpt_NET0\AdOpt_results\MY\EmissionLimit Brownfield\Zeeland\20250428122724_2040_minC_DD10-1\optimization_results.h5"
$ws.Range("G5").Value = "Z:\AdOpt_NET0\AdOpt_results\MY\EmissionLimit Brownfield\Zeeland\20250501121718_2050_minC_DD10-1\optimization_results.h5"

# Rows 6-95: refreshed numeric results from the new optimization runs
$ws.Range("E6").Value = [double]"6377673518.162631"
$ws.Range("F6").Value = [double]"5450814552.258259"
$ws.Range("G6").Value = [double]"3760890499.657729"
$ws.Range("F7").Value = [double]"537546119.7996608"
$ws.Range("G7").Value = [double]"875326557.1320901"
$ws.Range("E8").Value = [double]"6377673518.162631"
$ws.Range("F8").Value = [double]"5988360672.05792"
$ws.Range("G8").Value = [double]"4636217056.78982"
$ws.Range("G9").Value = [double]"164647051272.1071"
$ws.Range("E10").Value = [double]"1308190.96367628"
$ws.Range("F10").Value = [double]"654095.4818381065"
$ws.Range("G10").Value = [double]"0"
$ws.Range("E11").Value = [double]"894.6947646023427"
$ws.Range("E12").Value = [double]"43.368"
$ws.Range("F12").Value = [double]"1.496635535656232"
$ws.Range("E13").Value = [double]"383.0596234757299"
$ws.Range("F13").Value = [double]"6.767893059495433"
$ws.Range("G13").Value = [double]"155.5646101381337"
$ws.Range("E14").Value = [double]"502.5368929636838"
$ws.Range("F15").Value = [double]"210.605109753033"
$ws.Range("G15").Value = [double]"142.9572210309153"
$ws.Range("E20").Value = [double]"894.9482235336452"
$ws.Range("F20").Value = [double]"0"
$ws.Range("E21").Value = [double]"1248"
$ws.Range("F21").Value = [double]"43.06864850809301"
$ws.Range("E22").Value = [double]"0"
$ws.Range("F22").Value = [double]"43.06864850809301"
$ws.Range("E23").Value = [double]"469"
$ws.Range("E24").Value = [double]"689.4300000000668"
$ws.Range("E27").Value = [double]"88.3200931315462"
$ws.Range("E31").Value = [double]"5113.471999999979"
$ws.Range("F31").Value = [double]"10714.65279238021"
$ws.Range("E32").Value = [double]"1581.396712910418"
$ws.Range("E33").Value = [double]"182.3734825774716"
$ws.Range("F33").Value = [double]"1331.889642902058"
$ws.Range("E34").Value = [double]"0"
$ws.Range("F34").Value = [double]"5724.855129380961"
$ws.Range("E35").Value = [double]"1715.087259231281"
$ws.Range("F35").Value = [double]"13841.18942943548"
$ws.Range("E36").Value = [double]"0"
$ws.Range("F36").Value = [double]"3790.473682560014"
$ws.Range("E37").Value = [double]"167.3553178007917"
$ws.Range("F37").Value = [double]"0"
$ws.Range("E38").Value = [double]"1350.322154607045"
$ws.Range("E42").Value = [double]"259.6151598489242"
$ws.Range("F42").Value = [double]"242.9559191542486"
$ws.Range("G42").Value = [double]"407.3118751406292"
$ws.Range("E43").Value = [double]"0"
$ws.Range("F43").Value = [double]"2.816591404553037e-11"
$ws.Range("F49").Value = [double]"469.0000000000455"
$ws.Range("F57").Value = [double]"3250.00000000464"
$ws.Range("G57").Value = [double]"4600"
$ws.Range("E66").Value = [double]"145.8352466307604"
$ws.Range("F66").Value = [double]"145.8352466301818"
$ws.Range("G66").Value = [double]"145.8352466301818"
$ws.Range("E69").Value = [double]"1350.322154607045"
$ws.Range("F69").Value = [double]"1210.649120408641"
$ws.Range("G69").Value = [double]"894.9482235336452"
$ws.Range("F80").Value = [double]"6.084022174945858e-14"
$ws.Range("E84").Value = [double]"186.5220559262779"
$ws.Range("F84").Value = [double]"188.3584277279992"
$ws.Range("G84").Value = [double]"186.7105215149341"
$ws.Range("E85").Value = [double]"88.32009313296318"
$ws.Range("F95").Value = [double]"146.2075946034622"
$ws.Range("G95").Value = [double]"107.9031888874626"
$ws.Range("E98").Value = [double]"1248"
$ws.Range("F98").Value = [double]"1248"
$ws.Range("G98").Value = [double]"1291.068648508093"
$ws.Range("F101").Value = [double]"894.6947646023427"
$ws.Range("G101").Value = [double]"894.6947646023427"
$ws.Range("F102").Value = [double]"43.368"
$ws.Range("G102").Value = [double]"44.86463553565623"
$ws.Range("F103").Value = [double]"383.0596234757299"
$ws.Range("G103").Value = [double]"389.8275165352254"
$ws.Range("F104").Value = [double]"502.5368929636838"
$ws.Range("G104").Value = [double]"502.5368929636838"
$ws.Range("F105").Value = [double]"894.9482235336452"
$ws.Range("G105").Value = [double]"894.9482235336452"
$ws.Range("F106").Value = [double]"1248"
$ws.Range("G106").Value = [double]"1291.068648508093"
$ws.Range("F107").Value = [double]"469"
$ws.Range("G107").Value = [double]"469"
$ws.Range("F108").Value = [double]"689.4300000000668"
$ws.Range("G108").Value = [double]"689.4300000000668"
$ws.Range("F109").Value = [double]"88.3200931315462"
$ws.Range("G109").Value = [double]"88.3200931315462"
$ws.Range("F111").Value = [double]"5113.471999999979"
$ws.Range("G111").Value = [double]"15828.12479238019"
$ws.Range("F112").Value = [double]"1581.396712910418"
$ws.Range("G112").Value = [double]"5646.396712910419"
$ws.Range("F113").Value = [double]"182.3734825774716"
$ws.Range("G113").Value = [double]"1514.26312547953"

# Rows 114-120: "size_Storage_Ethylene_existing" row moved down (was row 114, now row 119);
# labels in rows 114-118 shift up one slot, values refreshed throughout
$ws.Range("A114").Value = "size_Storage_H2_existing"
$ws.Range("F114").Value = [double]"1715.087259231281"
$ws.Range("G114").Value = [double]"15556.27668866676"
$ws.Range("A115").Value = "size_WGS_m_existing"
$ws.Range("F115").Value = [double]"167.3553178007917"
$ws.Range("G115").Value = [double]"167.3553178007917"
$ws.Range("A116").Value = "size_feedgas_mixer_existing"
$ws.Range("F116").Value = [double]"1350.322154607045"
$ws.Range("G116").Value = [double]"1350.322154607045"
$ws.Range("A117").Value = "size_CO2_mixer_existing"
$ws.Range("F117").Value = ""
$ws.Range("G117").Value = [double]"210.605109753033"
$ws.Range("A118").Value = "size_CO2electrolysis_existing"
$ws.Range("G118").Value = [double]"146.2075946034622"
$ws.Range("A119").Value = "size_Storage_Ethylene_existing"
$ws.Range("G119").Value = [double]"5724.855129380961"
$ws.Range("G120").Value = [double]"3790.473682560014"
